$wb = $excel.ActiveWorkbook

# --- Update note for ELECTRIC_TIME_STEP on sheet CONDUCTOR_input (D23) ---
$wsInput = $wb.Worksheets.Item("CONDUCTOR_input")
$wsInput.Range("D23").Value2 = "time step value for the electric transient solution. If None, uses the default value of dt_th/10, being dt_th the thermal time step."

# --- Update note for ELECTRIC_SOLVER on sheet CONDUCTOR_operation (D10) ---
$wsOperation = $wb.Worksheets.Item("CONDUCTOR_operation")
$wsOperation.Range("D10").Value2 = "Flag to select the solver for the electric module. Possible values: 0= steady state; 1 = transient. Defaults to 1. N.B. At the time being the steady solution is not availabye, althoug it works, to be consistent with the thermal hydrauilc solution."
